# Parameter Explanations.xlsx - commit: "Updated init skip steps slider"
#
# Adds three new parameter rows (create_mp4 / mp4_fps / save_all_images)
# beneath the existing "Init Steps Skipped" row, and updates the default
# "Diffusion Steps" value from 300 to 250.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update default value for the "Diffusion Steps" parameter (row 5) ---
$ws.Range("D5").Value = 250

# --- New rows 23-25: fill parameter_name (A) and label (B) first ---
$ws.Range("A23").Value = "create_mp4"
$ws.Range("B23").Value = "Video Generation"
$ws.Range("A24").Value = "mp4_fps"
$ws.Range("B24").Value = "Video FPS"
$ws.Range("A25").Value = "save_all_images"
$ws.Range("B25").Value = "Intermediate Images"

# --- Fill type (C) and default_value (D) columns ---
$ws.Range("C23").Value = "boolean"
$ws.Range("D23").Value = $true
$ws.Range("C24").Value = "number"
$ws.Range("D24").Value = 30
$ws.Range("C25").Value = "boolean"
$ws.Range("D25").Value = $false

# --- Fill description (E) column: save_all_images, create_mp4, mp4_fps ---
$ws.Range("E25").Value = "Whether or not to save image files for all of the diffusion steps, or just the final ones. If unchecked, intermediate images will delete once diffusion is complete."
$ws.Range("E23").Value = "Whether or not you want to save a video showing the start-to-finish diffusion process."
$ws.Range("E24").Value = "How many frames per second the video is."

# --- Update the active selection to match where the author ended up editing ---
$ws.Range("D28").Select()
